$d = $word.ActiveDocument

# --- Title paragraph: "ProjectPlan" -> "GGD" (drop leading space on next run) ---
$null = $d.Content.Find.Execute("ProjectPlan", $true, $false, $false, $false, $false, $true, 1, $false, "GGD ", 2)
$null = $d.Content.Find.Execute(" Meesterproef Niels Bron ", $true, $false, $false, $false, $false, $true, 1, $false, "Meesterproef Niels Bron ", 2)
